$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: simple integers 1-4 ---
$ws.Range("A1").Value = 1
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 3
$ws.Range("A4").Value = 4

# --- Column B: date values formatted as short dates ---
$ws.Range("B1").Value = 36785
$ws.Range("B1").NumberFormat = "mm-dd-yy"
$ws.Range("B1").Copy()
$ws.Range("B2:B4").PasteSpecial(-4122)
$ws.Range("B2").Value = 27523
$ws.Range("B3").Value = 25809
$ws.Range("B4").Value = 41436

# --- Column C: names (shared strings) ---
# Entered out of row order so the shared-string table matches the
# original authoring order (Maria E, Andres, Alejandro, Esteban).
$ws.Range("C2").Value = "Maria E"
$ws.Range("C1").Value = "Andres"
$ws.Range("C3").Value = "Alejandro"
$ws.Range("C4").Value = "Esteban"

# --- Column D: integers ---
$ws.Range("D1").Value = 17
$ws.Range("D2").Value = 43
$ws.Range("D3").Value = 48
$ws.Range("D4").Value = 4

# --- Row 5: empty, underlined cell (case for a formatted-but-blank cell) ---
$ws.Range("B5").Font.Underline = 1

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Selection moves to B1 ---
$ws.Range("B1").Select()
